$d = $word.ActiveDocument

function Find-ParagraphIndex($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text.Contains($needle)) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1. Remove the "Meta description" paragraph that currently follows the main
#    H1 title ("Play Bier Haus Oktoberfest Slot for Free - Review and
#    Features" / "Meta description" / ": Play the Bier Haus Oktoberfest ...").
# ---------------------------------------------------------------------------
$metaIdx = Find-ParagraphIndex $d "Meta description"
$metaPara = $d.Paragraphs.Item($metaIdx)
$d.Range($metaPara.Range.Start, $metaPara.Range.End).Delete()

# ---------------------------------------------------------------------------
# 2. The final paragraph of the document currently holds the DALLE image
#    prompt. We need to:
#      a) insert a new paragraph right before it containing a bold repeat of
#         the page title ("Play Bier Haus Oktoberfest Slot for Free -
#         Review and Features")
#      b) replace the DALLE-prompt paragraph's own text with the old meta
#         description sentence, keeping the paragraph's italic formatting.
# ---------------------------------------------------------------------------

# Use an existing plain-formatted paragraph (no bold/italic) as a copy
# template, so the duplicated paragraph does not inherit unwanted character
# formatting (e.g. list-bullet style or italics from neighbouring runs).
# The paragraph that begins "Are you ready to raise your stein" is a normal,
# unformatted body paragraph, so it is a safe structural template.
$templateIdx = Find-ParagraphIndex $d "Are you ready to raise your stein"
$template = $d.Paragraphs.Item($templateIdx)
$d.Range($template.Range.Start, $template.Range.End).Copy()

$dalleIdx = Find-ParagraphIndex $d "Prompt for DALLE"
$dallePara = $d.Paragraphs.Item($dalleIdx)
$insertPoint = $d.Range($dallePara.Range.Start, $dallePara.Range.Start)
$insertPoint.Paste()

# The pasted paragraph now sits at $dalleIdx; the original DALLE-prompt
# paragraph has been pushed down to $dalleIdx + 1.
$newPara = $d.Paragraphs.Item($dalleIdx)
$newRange = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$newRange.Text = "Play Bier Haus Oktoberfest Slot for Free - Review and Features"

$newPara = $d.Paragraphs.Item($dalleIdx)
$newRange = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$newRange.Font.Bold = 1

$dallePara = $d.Paragraphs.Item($dalleIdx + 1)
$dalleRange = $d.Range($dallePara.Range.Start, $dallePara.Range.End - 1)
$dalleRange.Text = "Play the Bier Haus Oktoberfest online slot game for free and discover its multiple bonuses and potential big payouts. Review of the game features and graphics."

Write-Host "Final paragraph count:" $d.Paragraphs.Count
